$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value2 = 76.22221999999999
$ws.Range("I9").Value2 = 58
$ws.Range("J9").Value2 = 99
$ws.Range("K9").Value2 = 58
$ws.Range("L9").Value2 = 99
$ws.Range("M9").Value2 = 111
$ws.Range("N9").Value2 = -437

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 1782.6
$ws.Range("I12").Value2 = 1728.25
$ws.Range("J12").Value2 = 2000
$ws.Range("K12").Value2 = 1728.25
$ws.Range("L12").Value2 = 2000
$ws.Range("M12").Value2 = -1558.25
$ws.Range("N12").Value2 = -2340

# ALC row 16
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value2 = 2000
$ws.Range("I16").Value2 = 2000
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 2000
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -1770

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 4747.25
$ws.Range("I137").Value2 = 5331
$ws.Range("J137").Value2 = 2996
$ws.Range("K137").Value2 = 15993
$ws.Range("L137").Value2 = 8988
$ws.Range("M137").Value2 = -13443
$ws.Range("N137").Value2 = -14088

# ARM row 25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value2 = 2000
$ws.Range("I25").Value2 = 2000
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 2000
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = -1598

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 3614.0588
$ws.Range("I74").Value2 = 3652.4375
$ws.Range("J74").Value2 = 3000
$ws.Range("K74").Value2 = 3652.4375
$ws.Range("L74").Value2 = 3000
$ws.Range("M74").Value2 = -2778.4375
$ws.Range("N74").Value2 = -4748

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 3614.0588
$ws.Range("I77").Value2 = 3652.4375
$ws.Range("J77").Value2 = 3000
$ws.Range("K77").Value2 = 18262.1875
$ws.Range("L77").Value2 = 15000
$ws.Range("M77").Value2 = -13894.1875
$ws.Range("N77").Value2 = -23736

# ARM row 98
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value2 = 15333
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 15333
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 15333
$ws.Range("N98").Value2 = -21323

# ARM row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value2 = 283351
$ws.Range("I101").Value2 = 0
$ws.Range("J101").Value2 = 283351
$ws.Range("K101").Value2 = 0
$ws.Range("L101").Value2 = 283351
$ws.Range("N101").Value2 = -289841

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value2 = 647
$ws.Range("I110").Value2 = 641.25
$ws.Range("J110").Value2 = 670
$ws.Range("K110").Value2 = 641.25
$ws.Range("L110").Value2 = 670
$ws.Range("M110").Value2 = 1403.75
$ws.Range("N110").Value2 = -4760

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 7666.6665
$ws.Range("I20").Value2 = 7500
$ws.Range("J20").Value2 = 8000
$ws.Range("K20").Value2 = 7500
$ws.Range("L20").Value2 = 8000
$ws.Range("M20").Value2 = -7253
$ws.Range("N20").Value2 = -8494

# BSM row 37
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value2 = 496
$ws.Range("I37").Value2 = 496
$ws.Range("J37").Value2 = 0
$ws.Range("K37").Value2 = 496
$ws.Range("L37").Value2 = 0
$ws.Range("M37").Value2 = -359

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 2265.75
$ws.Range("I107").Value2 = 2265.75
$ws.Range("J107").Value2 = 0
$ws.Range("K107").Value2 = 2265.75
$ws.Range("L107").Value2 = 0
$ws.Range("M107").Value2 = -345.75

# BSM row 117
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value2 = 71500
$ws.Range("I117").Value2 = 0
$ws.Range("J117").Value2 = 71500
$ws.Range("K117").Value2 = 0
$ws.Range("L117").Value2 = 71500
$ws.Range("N117").Value2 = -80678

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 872.4375
$ws.Range("I16").Value2 = 872.4375
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 872.4375
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -585.4375
$ws.Range("N16").ClearContents()

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 1750
$ws.Range("I22").Value2 = 1500
$ws.Range("J22").Value2 = 2000
$ws.Range("K22").Value2 = 1500
$ws.Range("L22").Value2 = 2000
$ws.Range("M22").Value2 = -1150
$ws.Range("N22").Value2 = -2700

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 0
$ws.Range("I58").Value2 = 0
$ws.Range("J58").Value2 = 0
$ws.Range("K58").Value2 = 0
$ws.Range("L58").Value2 = 0
$ws.Range("N58").ClearContents()

# CRP row 106
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value2 = 82798.8
$ws.Range("I106").Value2 = 0
$ws.Range("J106").Value2 = 82798.8
$ws.Range("K106").Value2 = 0
$ws.Range("L106").Value2 = 82798.8
$ws.Range("N106").Value2 = -85322.8

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 629.1429000000001
$ws.Range("I107").Value2 = 540.8
$ws.Range("J107").Value2 = 850
$ws.Range("K107").Value2 = 540.8
$ws.Range("L107").Value2 = 850
$ws.Range("M107").Value2 = 1379.2
$ws.Range("N107").Value2 = -4690

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value2 = 872.4375
$ws.Range("I113").Value2 = 872.4375
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 872.4375
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = 1297.5625
$ws.Range("N113").ClearContents()

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 1940.5
$ws.Range("I132").Value2 = 1825.2858
$ws.Range("J132").Value2 = 2209.3333
$ws.Range("K132").Value2 = 5475.857400000001
$ws.Range("L132").Value2 = 6627.999899999999
$ws.Range("M132").Value2 = -2945.857400000001
$ws.Range("N132").Value2 = -11687.9999

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 0
$ws.Range("I136").Value2 = 0
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 0
$ws.Range("L136").Value2 = 0
$ws.Range("N136").ClearContents()

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 20.714285
$ws.Range("I2").Value2 = 20.333334
$ws.Range("J2").Value2 = 21
$ws.Range("K2").Value2 = 122.000004
$ws.Range("L2").Value2 = 126
$ws.Range("M2").Value2 = -9.000004000000004
$ws.Range("N2").Value2 = -352

# CUL row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value2 = 3608.875
$ws.Range("I14").Value2 = 3608.875
$ws.Range("J14").Value2 = 0
$ws.Range("K14").Value2 = 10826.625
$ws.Range("L14").Value2 = 0
$ws.Range("M14").Value2 = -10653.625

# CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value2 = 637.1429000000001
$ws.Range("I18").Value2 = 714.2
$ws.Range("J18").Value2 = 444.5
$ws.Range("K18").Value2 = 2142.6
$ws.Range("L18").Value2 = 1333.5
$ws.Range("M18").Value2 = -1973.6
$ws.Range("N18").Value2 = -1671.5

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value2 = 1245
$ws.Range("I92").Value2 = 993.3333
$ws.Range("J92").Value2 = 2000
$ws.Range("K92").Value2 = 2979.9999
$ws.Range("L92").Value2 = 6000
$ws.Range("M92").Value2 = -1731.9999
$ws.Range("N92").Value2 = -8496

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 4188000.5
$ws.Range("I11").Value2 = 7571251
$ws.Range("J11").Value2 = 321428.56
$ws.Range("K11").Value2 = 7571251
$ws.Range("L11").Value2 = 321428.56
$ws.Range("M11").Value2 = -7571112
$ws.Range("N11").Value2 = -321706.56

# GSM row 21
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value2 = 966666.7
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value2 = 966666.7
$ws.Range("K21").Value2 = 0
$ws.Range("L21").Value2 = 966666.7
$ws.Range("N21").Value2 = -967012.7

# GSM row 30
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value2 = 966666.7
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 966666.7
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 966666.7
$ws.Range("N30").Value2 = -966876.7

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 2000
$ws.Range("I70").Value2 = 2000
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 2000
$ws.Range("L70").Value2 = 0
$ws.Range("M70").Value2 = -1730

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value2 = 2000
$ws.Range("I73").Value2 = 2000
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 2000
$ws.Range("L73").Value2 = 0
$ws.Range("M73").Value2 = -1064

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value2 = 722.125
$ws.Range("I107").Value2 = 883.3333
$ws.Range("J107").Value2 = 238.5
$ws.Range("K107").Value2 = 883.3333
$ws.Range("L107").Value2 = 238.5
$ws.Range("M107").Value2 = 1036.6667
$ws.Range("N107").Value2 = -4078.5

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 23872.166
$ws.Range("I40").Value2 = 14646.1
$ws.Range("J40").Value2 = 70002.5
$ws.Range("K40").Value2 = 14646.1
$ws.Range("L40").Value2 = 70002.5
$ws.Range("M40").Value2 = -14510.1
$ws.Range("N40").Value2 = -70274.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 335132.84
$ws.Range("I46").Value2 = 1000400
$ws.Range("J46").Value2 = 2499.25
$ws.Range("K46").Value2 = 1000400
$ws.Range("L46").Value2 = 2499.25
$ws.Range("M46").Value2 = -1000212
$ws.Range("N46").Value2 = -2875.25

# LTW row 104
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value2 = 28498.5
$ws.Range("I104").Value2 = 0
$ws.Range("J104").Value2 = 28498.5
$ws.Range("K104").Value2 = 0
$ws.Range("L104").Value2 = 28498.5
$ws.Range("N104").Value2 = -35486.5

# LTW row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value2 = 10000
$ws.Range("I106").Value2 = 0
$ws.Range("J106").Value2 = 10000
$ws.Range("K106").Value2 = 0
$ws.Range("L106").Value2 = 10000
$ws.Range("N106").Value2 = -12524

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 498
$ws.Range("I107").Value2 = 498
$ws.Range("J107").Value2 = 0
$ws.Range("K107").Value2 = 1494
$ws.Range("L107").Value2 = 0
$ws.Range("M107").Value2 = 426
$ws.Range("N107").ClearContents()

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value2 = 11701.2
$ws.Range("I113").Value2 = 876.5
$ws.Range("J113").Value2 = 55000
$ws.Range("K113").Value2 = 2629.5
$ws.Range("L113").Value2 = 165000
$ws.Range("M113").Value2 = -459.5
$ws.Range("N113").Value2 = -169340

# WVR row 116
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value2 = 51249.75
$ws.Range("I116").Value2 = 0
$ws.Range("J116").Value2 = 51249.75
$ws.Range("K116").Value2 = 0
$ws.Range("L116").Value2 = 51249.75
$ws.Range("N116").Value2 = -60427.75

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 280
$ws.Range("I122").Value2 = 555
$ws.Range("J122").Value2 = 5
$ws.Range("K122").Value2 = 1665
$ws.Range("L122").Value2 = 15
$ws.Range("M122").Value2 = 785
$ws.Range("N122").Value2 = -4915

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 12999.3
$ws.Range("I136").Value2 = 13332.333
$ws.Range("J136").Value2 = 12499.75
$ws.Range("K136").Value2 = 39996.999
$ws.Range("L136").Value2 = 37499.25
$ws.Range("M136").Value2 = -37446.999
$ws.Range("N136").Value2 = -42599.25
